# Change names from *img to img*
$wb = $excel.ActiveWorkbook

# Rename the "*img" sheets to "img*" (move the "img" prefix to a suffix
# position, preserving the trailing letter order: himg->imgh, timg->imgt,
# simg->imgs, gimg->imgg, wimg->imgw, bimg->imgb, eimg->imge).
$renames = @{
    "himg" = "imgh"
    "timg" = "imgt"
    "simg" = "imgs"
    "gimg" = "imgg"
    "wimg" = "imgw"
    "bimg" = "imgb"
    "eimg" = "imge"
}

foreach ($oldName in $renames.Keys) {
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $renames[$oldName]
}

# The last sheet (now named "imge") becomes the active/selected sheet.
$wb.Worksheets.Item("imge").Activate()
